$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ORDER")

# Header title: clarify the CUS_ID/ORDER_ID "register" key is now a random suffix,
# not tied to ORDER_ID (fixes the broken package-id calc referenced in the commit).
$ws.Range("C1").Value = "REGISTER + CUS_ID+RANDOM(6)"

# Existing rows 3-4 get the corrected, randomised order-id text.
$ws.Range("C3").Value = "20170420-1-111"
$ws.Range("C4").Value = "20170420-1-111"

# Duplicate row 4's layout/format down through row 8 so the new rows inherit
# the same styles (date format on D:F, centered style on G) as row 4.
$ws.Range("A4:G4").Copy() | Out-Null
$ws.Range("A5:G8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($i = 5; $i -le 8; $i++) {
    $b = $i - 2
    $ws.Cells.Item($i, 1).Value = 2
    $ws.Cells.Item($i, 2).Value = $b
    $ws.Cells.Item($i, 3).Value = "20170420-1-111"
    $ws.Cells.Item($i, 4).Value = 43332
    $ws.Cells.Item($i, 7).Value = 1
}

# Re-assert the date formulas across the whole block so Excel stores them as
# one shared formula (E3:E8 / F3:F8) like the original two-row block did.
$ws.Range("E3:E8").Formula = "=D3+3"
$ws.Range("F3:F8").Formula = "=D3+366"

# Column C needs to widen a bit to fit the new header text.
$ws.Columns.Item(3).ColumnWidth = 28.9

# Leave the selection where the user last worked (newly-added rows' helper column).
$ws.Activate() | Out-Null
$ws.Range("I5:I8").Select() | Out-Null

# The workbook window had been minimized when this was last saved.
$win = $wb.Windows.Item(1)
$win.WindowState = -4140
